# Lisää pohjavesialueiden siirrolle ja toteumatyökalulle omat käyttöoikeudet
#
# On the "Oikeudet" sheet, two new rows are inserted right after the existing
# "Järjestelmäasetukset" row (old row 111) and before the "Testaus" row (old
# row 112), pushing everything from the old row 112 onward down by two rows:
#   new row 112: Hallinta / Pohjavesialueidensiirto / ... / R*,W*
#   new row 113: Hallinta / Toteumatyokalu          / ... / R*,W*
# (old row 112 "Testaus" -> new row 114, old row 113 "Laadunseuranta" -> new row 115)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oikeudet")
$ws.Activate()

# --- 1. Insert two new rows before the old row 112 ("Testaus") ---------------
$ws.Range("A112:A113").EntireRow.Insert()

# Copy the formatting of row 111 ("Järjestelmäasetukset", which keeps the
# "Hallinta" section look) onto the two freshly-inserted rows so the new rows
# match the surrounding "Hallinta" block styling instead of Excel's bare
# insert-row default.
$ws.Range("A111:Y111").Copy()
$ws.Range("A112:Y113").PasteSpecial(-4122)
$ws.Rows.Item(112).RowHeight = 15.75
$ws.Rows.Item(113).RowHeight = 15.75
$excel.CutCopyMode = 0

# --- 2. Fill in the new rows' content ----------------------------------------
$ws.Range("A112").Value = "Hallinta"
$ws.Range("B112").Value = "Pohjavesialueidensiirto"
$ws.Range("D112").Value = "R*,W*"

$ws.Range("A113").Value = "Hallinta"
$ws.Range("B113").Value = "Toteumatyokalu"
$ws.Range("D113").Value = "R*,W*"

# --- 3. Grow the filter range / named ranges from row 113 to row 115 --------
$names = $wb.Names
for ($i = 1; $i -le $names.Count(); $i++) {
  $n = $names.Item($i)
  $old = $n.RefersTo()
  if ($old.IndexOf('Oikeudet!$A$5:$Y$113') -ge 0) {
    $n.RefersTo = $old.Replace('Oikeudet!$A$5:$Y$113', 'Oikeudet!$A$5:$Y$115')
  }
}

# --- 4. Move the selection to the newly-edited cell, matching the author's
#        last touched cell in the source edit ---------------------------------
$ws.Range("B113").Select()
